$p = $ppt.ActivePresentation
Write-Output "HasHandoutMaster: $($p.HasHandoutMaster)"
try {
  $hm = $p.HandoutMaster
  Write-Output "HandoutMaster: $hm"
  Write-Output "HandoutMaster.Name: $($hm.Name)"
} catch { Write-Output "err: $_" }
